$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '51.854.23'
$ws.Range('E2').NumberFormat = "@"
$ws.Range('E2').Value = '  +0.20%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.813.86'
$ws.Range('E3').NumberFormat = "@"
$ws.Range('E3').Value = '  +1.32%  '
$ws.Range('E4').NumberFormat = "@"
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '356.91'
$ws.Range('E5').NumberFormat = "@"
$ws.Range('E5').Value = '  +0.20%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '110.12'
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.558'
$ws.Range('E7').NumberFormat = "@"
$ws.Range('E7').Value = '  +0.57%  '
$ws.Range('E8').NumberFormat = "@"
$ws.Range('E8').Value = '  +0.06%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.634'
$ws.Range('E9').NumberFormat = "@"
$ws.Range('E9').Value = '  +8.39%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '40.19'
$ws.Range('E10').NumberFormat = "@"
$ws.Range('E10').Value = '  +1.33%  '
$ws.Range('E11').NumberFormat = "@"
$ws.Range('E11').Value = '  +0.17%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.0838'
$ws.Range('E12').NumberFormat = "@"
$ws.Range('E12').Value = '  -0.68%  '
$ws.Range('E13').NumberFormat = "@"
$ws.Range('E13').Value = '  +3.14%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '7.81'
$ws.Range('E14').NumberFormat = "@"
$ws.Range('E14').Value = '  +2.63%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '3.261.35'
$ws.Range('E15').NumberFormat = "@"
$ws.Range('E15').Value = '  +1.55%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '2.820.34'
$ws.Range('E16').NumberFormat = "@"
$ws.Range('E16').Value = '  +1.00%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.944'
$ws.Range('E17').NumberFormat = "@"
$ws.Range('E17').Value = '  +1.26%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '51.883.58'
$ws.Range('E18').NumberFormat = "@"
$ws.Range('E18').Value = '  +0.42%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '7.68'
$ws.Range('E19').NumberFormat = "@"
$ws.Range('E19').Value = '  +2.94%  '
$ws.Range('E20').NumberFormat = "@"
$ws.Range('E20').Value = '  +3.23%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '13.67'
$ws.Range('E21').NumberFormat = "@"
$ws.Range('E21').Value = '  +4.18%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.0₃0979'
$ws.Range('E22').NumberFormat = "@"
$ws.Range('E22').Value = '  +1.13%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '70.47'
$ws.Range('E23').NumberFormat = "@"
$ws.Range('E23').Value = '  +0.50%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '268.56'
$ws.Range('E24').NumberFormat = "@"
$ws.Range('E24').Value = '  -0.05%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.76'
$ws.Range('E25').NumberFormat = "@"
$ws.Range('E25').Value = '  +1.28%  '
$ws.Range('B26').Value = 'Dai'
$ws.Range('C26').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '1.00'
$ws.Range('E26').NumberFormat = "@"
$ws.Range('E26').Value = '  +0.07%  '
$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('C27').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '26.18'
$ws.Range('E27').NumberFormat = "@"
$ws.Range('E27').Value = '  -0.77%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '0.164'
$ws.Range('E28').NumberFormat = "@"
$ws.Range('E28').Value = '  +0.98%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '10.42'
$ws.Range('E29').NumberFormat = "@"
$ws.Range('E29').Value = '  +1.94%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '38.11'
$ws.Range('E30').NumberFormat = "@"
$ws.Range('E30').Value = '  +9.29%  '
$ws.Range('E31').NumberFormat = "@"
$ws.Range('E31').Value = '  +1.01%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '6.19'
$ws.Range('E32').NumberFormat = "@"
$ws.Range('E32').Value = '  -1.43%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '52.13'
$ws.Range('E33').NumberFormat = "@"
$ws.Range('E33').Value = '  +1.01%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '5.66'
$ws.Range('E34').NumberFormat = "@"
$ws.Range('E34').Value = '  +10.59%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.0447'
$ws.Range('E35').NumberFormat = "@"
$ws.Range('E35').Value = '  -0.42%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.0870'
$ws.Range('E36').NumberFormat = "@"
$ws.Range('E36').Value = '  +3.80%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '1.00'
$ws.Range('E37').NumberFormat = "@"
$ws.Range('E37').Value = '  +0.08%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '18.83'
$ws.Range('E38').NumberFormat = "@"
$ws.Range('E38').Value = '  +0.86%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '2.01'
$ws.Range('E39').NumberFormat = "@"
$ws.Range('E39').Value = '  +2.68%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '3.15'
$ws.Range('E40').NumberFormat = "@"
$ws.Range('E40').Value = '  +0.40%  '
$ws.Range('E41').NumberFormat = "@"
$ws.Range('E41').Value = '  +1.04%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '2.50'
$ws.Range('E42').NumberFormat = "@"
$ws.Range('E42').Value = '  -1.52%  '
$ws.Range('B43').Value = 'Monero'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '119.52'
$ws.Range('E43').NumberFormat = "@"
$ws.Range('E43').Value = '  +0.31%  '
$ws.Range('E44').NumberFormat = "@"
$ws.Range('E44').Value = '  -1.00%  '
$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '21.95'
$ws.Range('E45').NumberFormat = "@"
$ws.Range('E45').Value = '  +1.70%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '2.48'
$ws.Range('E46').NumberFormat = "@"
$ws.Range('E46').Value = '  +8.76%  '
$ws.Range('B47').Value = 'Maker'
$ws.Range('C47').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '2.109.83'
$ws.Range('E47').NumberFormat = "@"
$ws.Range('E47').Value = '  +1.29%  '
$ws.Range('B48').Value = 'NEARProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '3.39'
$ws.Range('E48').NumberFormat = "@"
$ws.Range('E48').Value = '  +3.56%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.930'
$ws.Range('E49').NumberFormat = "@"
$ws.Range('E49').Value = '  -1.44%  '
$ws.Range('E50').NumberFormat = "@"
$ws.Range('E50').Value = '  +10.00%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '5.45'
$ws.Range('E51').NumberFormat = "@"
$ws.Range('E51').Value = '  -2.47%  '
